$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 13500442  # H6 was 437.69232
$ws.Cells.Item(6, 9).Value = 16200340  # I6 was 424.2
$ws.Cells.Item(6, 10).Value = 950  # J6 was 482.66666
$ws.Cells.Item(6, 11).Value = 48601020  # K6 was 1272.6
$ws.Cells.Item(6, 12).Value = 2850  # L6 was 1447.99998
$ws.Cells.Item(6, 13).Value = -48600908  # M6 was -1160.6
$ws.Cells.Item(6, 14).Value = -3074  # N6 was -1671.99998

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 633.3333  # H12 was 350
$ws.Cells.Item(12, 9).Value = 200  # I12 was 0
$ws.Cells.Item(12, 10).Value = 850  # J12 was 350
$ws.Cells.Item(12, 11).Value = 200  # K12 was 0
$ws.Cells.Item(12, 12).Value = 850  # L12 was 350
$ws.Cells.Item(12, 13).Value = -30  # M12 was None
$ws.Cells.Item(12, 14).Value = -1190  # N12 was -690

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1247.5  # H29 was 1158
$ws.Cells.Item(29, 9).Value = 1000  # I29 was 950
$ws.Cells.Item(29, 10).Value = 1990  # J29 was 1990
$ws.Cells.Item(29, 11).Value = 3000  # K29 was 2850
$ws.Cells.Item(29, 12).Value = 5970  # L29 was 5970
$ws.Cells.Item(29, 13).Value = -2719  # M29 was -2569
$ws.Cells.Item(29, 14).Value = -6532  # N29 was -6532

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 75.5  # H38 was 78.75
$ws.Cells.Item(38, 9).Value = 75.5  # I38 was 78.75
$ws.Cells.Item(38, 10).Value = 0  # J38 was 0
$ws.Cells.Item(38, 11).Value = 226.5  # K38 was 236.25
$ws.Cells.Item(38, 12).Value = 0  # L38 was 0
$ws.Cells.Item(38, 13).Value = 145.5  # M38 was 135.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 478.75  # H58 was 1829.5
$ws.Cells.Item(58, 9).Value = 478.75  # I58 was 583
$ws.Cells.Item(58, 10).Value = 0  # J58 was 3076
$ws.Cells.Item(58, 11).Value = 1436.25  # K58 was 1749
$ws.Cells.Item(58, 12).Value = 0  # L58 was 9228
$ws.Cells.Item(58, 13).ClearContents()  # M58 was -1599
$ws.Cells.Item(58, 14).Value = -1286.25  # N58 was -9528

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3673.0667  # H64 was 3656.75
$ws.Cells.Item(64, 9).Value = 3866.3333  # I64 was 3722.25
$ws.Cells.Item(64, 10).Value = 3651.5925  # J64 was 3645.8333
$ws.Cells.Item(64, 11).Value = 3866.3333  # K64 was 3722.25
$ws.Cells.Item(64, 12).Value = 3651.5925  # L64 was 3645.8333
$ws.Cells.Item(64, 13).Value = -3618.3333  # M64 was -3474.25
$ws.Cells.Item(64, 14).Value = -4147.592500000001  # N64 was -4141.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 3673.0667  # H67 was 3656.75
$ws.Cells.Item(67, 9).Value = 3866.3333  # I67 was 3722.25
$ws.Cells.Item(67, 10).Value = 3651.5925  # J67 was 3645.8333
$ws.Cells.Item(67, 11).Value = 3866.3333  # K67 was 3722.25
$ws.Cells.Item(67, 12).Value = 3651.5925  # L67 was 3645.8333
$ws.Cells.Item(67, 13).Value = -3008.3333  # M67 was -2864.25
$ws.Cells.Item(67, 14).Value = -5367.592500000001  # N67 was -5361.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 10).Value = 4333.3335  # J76 was 4342.857
$ws.Cells.Item(76, 11).Value = 4800  # K76 was 4800
$ws.Cells.Item(76, 12).Value = 4333.3335  # L76 was 4342.857
$ws.Cells.Item(76, 13).Value = -4485  # M76 was -4485
$ws.Cells.Item(76, 14).Value = -4963.3335  # N76 was -4972.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 10).Value = 4333.3335  # J79 was 4342.857
$ws.Cells.Item(79, 11).Value = 4800  # K79 was 4800
$ws.Cells.Item(79, 12).Value = 4333.3335  # L79 was 4342.857
$ws.Cells.Item(79, 13).Value = -3708  # M79 was -3708
$ws.Cells.Item(79, 14).Value = -6517.3335  # N79 was -6526.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 18713  # H82 was 16649.143
$ws.Cells.Item(82, 9).Value = 1617.3334  # I82 was 2308.8
$ws.Cells.Item(82, 10).Value = 70000  # J82 was 52500
$ws.Cells.Item(82, 11).Value = 4852.0002  # K82 was 6926.400000000001
$ws.Cells.Item(82, 12).Value = 210000  # L82 was 157500
$ws.Cells.Item(82, 13).Value = -4446.0002  # M82 was -6520.400000000001
$ws.Cells.Item(82, 14).Value = -210812  # N82 was -158312

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(85, 8).Value = 18713  # H85 was 16649.143
$ws.Cells.Item(85, 9).Value = 1617.3334  # I85 was 2308.8
$ws.Cells.Item(85, 10).Value = 70000  # J85 was 52500
$ws.Cells.Item(85, 11).Value = 4852.0002  # K85 was 6926.400000000001
$ws.Cells.Item(85, 12).Value = 210000  # L85 was 157500
$ws.Cells.Item(85, 13).Value = -3448.0002  # M85 was -5522.400000000001
$ws.Cells.Item(85, 14).Value = -212808  # N85 was -160308

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 2516.9167  # H88 was 2478.5715
$ws.Cells.Item(88, 9).Value = 1601  # I88 was 1833.3334
$ws.Cells.Item(88, 10).Value = 2822.2222  # J88 was 2962.5
$ws.Cells.Item(88, 11).Value = 1601  # K88 was 1833.3334
$ws.Cells.Item(88, 12).Value = 2822.2222  # L88 was 2962.5
$ws.Cells.Item(88, 13).Value = -1195  # M88 was -1427.3334
$ws.Cells.Item(88, 14).Value = -3634.2222  # N88 was -3774.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 2516.9167  # H91 was 2478.5715
$ws.Cells.Item(91, 9).Value = 1601  # I91 was 1833.3334
$ws.Cells.Item(91, 10).Value = 2822.2222  # J91 was 2962.5
$ws.Cells.Item(91, 11).Value = 1601  # K91 was 1833.3334
$ws.Cells.Item(91, 12).Value = 2822.2222  # L91 was 2962.5
$ws.Cells.Item(91, 13).Value = -197  # M91 was -429.3334
$ws.Cells.Item(91, 14).Value = -5630.2222  # N91 was -5770.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 10769.692  # H112 was 10064.714
$ws.Cells.Item(112, 9).Value = 875  # I112 was 875
$ws.Cells.Item(112, 10).Value = 15167.333  # J112 was 13740.6
$ws.Cells.Item(112, 11).Value = 2625  # K112 was 2625
$ws.Cells.Item(112, 12).Value = 45501.999  # L112 was 41221.8
$ws.Cells.Item(112, 13).Value = -1517  # M112 was -1517
$ws.Cells.Item(112, 14).Value = -47717.999  # N112 was -43437.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 2431.4285  # H118 was 2145.8333
$ws.Cells.Item(118, 9).Value = 466.66666  # I118 was 497.5
$ws.Cells.Item(118, 10).Value = 2967.2727  # J118 was 2970
$ws.Cells.Item(118, 11).Value = 1399.99998  # K118 was 1492.5
$ws.Cells.Item(118, 12).Value = 8901.8181  # L118 was 8910
$ws.Cells.Item(118, 13).Value = 257.0000199999999  # M118 was 164.5
$ws.Cells.Item(118, 14).Value = -12215.8181  # N118 was -12224

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(123, 8).Value = 65856.92  # H123 was 65832.63
$ws.Cells.Item(123, 9).Value = 0  # I123 was 0
$ws.Cells.Item(123, 10).Value = 65856.92  # J123 was 65832.63
$ws.Cells.Item(123, 11).Value = 0  # K123 was 0
$ws.Cells.Item(123, 12).Value = 65856.92  # L123 was 65832.63
$ws.Cells.Item(123, 14).Value = -75656.92  # N123 was -75632.63

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 0  # H124 was 60000
$ws.Cells.Item(124, 9).Value = 0  # I124 was 0
$ws.Cells.Item(124, 10).Value = 0  # J124 was 60000
$ws.Cells.Item(124, 11).Value = 0  # K124 was 0
$ws.Cells.Item(124, 12).ClearContents()  # L124 was 60000
$ws.Cells.Item(124, 14).Value = 0  # N124 was -69820

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value = 36312  # H126 was 40115.555
$ws.Cells.Item(126, 9).Value = 20000  # I126 was 0
$ws.Cells.Item(126, 10).Value = 40390  # J126 was 40115.555
$ws.Cells.Item(126, 11).Value = 20000  # K126 was 0
$ws.Cells.Item(126, 12).Value = 40390  # L126 was 40115.555
$ws.Cells.Item(126, 13).Value = -15060  # M126 was None
$ws.Cells.Item(126, 14).Value = -50270  # N126 was -49995.555

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1220.8667  # H129 was 1339.125
$ws.Cells.Item(129, 9).Value = 740  # I129 was 740
$ws.Cells.Item(129, 10).Value = 1341.0834  # J129 was 1698.6
$ws.Cells.Item(129, 11).Value = 2220  # K129 was 2220
$ws.Cells.Item(129, 12).Value = 4023.2502  # L129 was 5095.799999999999
$ws.Cells.Item(129, 13).Value = 2780  # M129 was 2780
$ws.Cells.Item(129, 14).Value = -14023.2502  # N129 was -15095.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4023.11  # H138 was 3962.94
$ws.Cells.Item(138, 9).Value = 3185.152  # I138 was 3145.0425
$ws.Cells.Item(138, 10).Value = 4736.926  # J138 was 4688.245
$ws.Cells.Item(138, 11).Value = 9555.456  # K138 was 9435.127500000001
$ws.Cells.Item(138, 12).Value = 14210.778  # L138 was 14064.735
$ws.Cells.Item(138, 13).Value = -4415.456  # M138 was -4295.127500000001
$ws.Cells.Item(138, 14).Value = -24490.778  # N138 was -24344.735

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 84817  # H32 was 32628.174
$ws.Cells.Item(32, 9).Value = 103489.14  # I32 was 25782.904
$ws.Cells.Item(32, 10).Value = 66144.86  # J32 was 104503.5
$ws.Cells.Item(32, 11).Value = 103489.14  # K32 was 25782.904
$ws.Cells.Item(32, 12).Value = 66144.86  # L32 was 104503.5
$ws.Cells.Item(32, 13).Value = -103202.14  # M32 was -25495.904
$ws.Cells.Item(32, 14).Value = -66718.86  # N32 was -105077.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2500  # H88 was 2669
$ws.Cells.Item(88, 9).Value = 0  # I88 was 0
$ws.Cells.Item(88, 10).Value = 2500  # J88 was 2669
$ws.Cells.Item(88, 11).Value = 0  # K88 was 0
$ws.Cells.Item(88, 12).Value = 2500  # L88 was 2669
$ws.Cells.Item(88, 14).Value = -3312  # N88 was -3481

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 2500  # H91 was 2669
$ws.Cells.Item(91, 9).Value = 0  # I91 was 0
$ws.Cells.Item(91, 10).Value = 2500  # J91 was 2669
$ws.Cells.Item(91, 11).Value = 0  # K91 was 0
$ws.Cells.Item(91, 12).Value = 2500  # L91 was 2669
$ws.Cells.Item(91, 14).Value = -5308  # N91 was -5477

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 549.3333  # H97 was 613.55554
$ws.Cells.Item(97, 9).Value = 505.5  # I97 was 613.55554
$ws.Cells.Item(97, 10).Value = 900  # J97 was 0
$ws.Cells.Item(97, 11).Value = 505.5  # K97 was 613.55554
$ws.Cells.Item(97, 12).Value = 900  # L97 was 0
$ws.Cells.Item(97, 13).Value = -9.5  # M97 was -117.55554
$ws.Cells.Item(97, 14).Value = -1892  # N97 was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(69, 8).Value = 0  # H69 was 32000
$ws.Cells.Item(69, 9).Value = 0  # I69 was 0
$ws.Cells.Item(69, 10).Value = 0  # J69 was 32000
$ws.Cells.Item(69, 11).Value = 0  # K69 was 0
$ws.Cells.Item(69, 12).ClearContents()  # L69 was 32000
$ws.Cells.Item(69, 14).Value = 0  # N69 was -33622

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(72, 8).Value = 0  # H72 was 32000
$ws.Cells.Item(72, 9).Value = 0  # I72 was 0
$ws.Cells.Item(72, 10).Value = 0  # J72 was 32000
$ws.Cells.Item(72, 11).Value = 0  # K72 was 0
$ws.Cells.Item(72, 12).ClearContents()  # L72 was 96000
$ws.Cells.Item(72, 14).Value = 0  # N72 was -104112

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 264512.75  # H86 was 667002.3
$ws.Cells.Item(86, 9).Value = 35099.332  # I86 was 1000
$ws.Cells.Item(86, 10).Value = 402160.8  # J86 was 1000003.5
$ws.Cells.Item(86, 11).Value = 35099.332  # K86 was 1000
$ws.Cells.Item(86, 12).Value = 402160.8  # L86 was 1000003.5
$ws.Cells.Item(86, 13).Value = -33976.332  # M86 was 123
$ws.Cells.Item(86, 14).Value = -404406.8  # N86 was -1002249.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 264512.75  # H89 was 667002.3
$ws.Cells.Item(89, 9).Value = 35099.332  # I89 was 1000
$ws.Cells.Item(89, 10).Value = 402160.8  # J89 was 1000003.5
$ws.Cells.Item(89, 11).Value = 175496.66  # K89 was 5000
$ws.Cells.Item(89, 12).Value = 2010804  # L89 was 5000017.5
$ws.Cells.Item(89, 13).Value = -169880.66  # M89 was 616
$ws.Cells.Item(89, 14).Value = -2022036  # N89 was -5011249.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 53745.26  # H94 was 46330.41
$ws.Cells.Item(94, 9).Value = 820  # I94 was 653.8333
$ws.Cells.Item(94, 10).Value = 168416.67  # J94 was 251875
$ws.Cells.Item(94, 11).Value = 820  # K94 was 653.8333
$ws.Cells.Item(94, 12).Value = 168416.67  # L94 was 251875
$ws.Cells.Item(94, 13).Value = -369  # M94 was -202.8333
$ws.Cells.Item(94, 14).Value = -169318.67  # N94 was -252777

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1622.7142  # H99 was 1693.6923
$ws.Cells.Item(99, 9).Value = 1091.8  # I99 was 1135.3334
$ws.Cells.Item(99, 10).Value = 2950  # J99 was 2950
$ws.Cells.Item(99, 11).Value = 1091.8  # K99 was 1135.3334
$ws.Cells.Item(99, 12).Value = 2950  # L99 was 2950
$ws.Cells.Item(99, 13).Value = 406.2  # M99 was 362.6666
$ws.Cells.Item(99, 14).Value = -5946  # N99 was -5946

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5366.5  # H31 was 5616.5
$ws.Cells.Item(31, 9).Value = 4179.5  # I31 was 4519.4287
$ws.Cells.Item(31, 10).Value = 6553.5  # J31 was 6469.778
$ws.Cells.Item(31, 11).Value = 4179.5  # K31 was 4519.4287
$ws.Cells.Item(31, 12).Value = 6553.5  # L31 was 6469.778
$ws.Cells.Item(31, 13).Value = -3884.5  # M31 was -4224.4287
$ws.Cells.Item(31, 14).Value = -7143.5  # N31 was -7059.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5366.5  # H34 was 5616.5
$ws.Cells.Item(34, 9).Value = 4179.5  # I34 was 4519.4287
$ws.Cells.Item(34, 10).Value = 6553.5  # J34 was 6469.778
$ws.Cells.Item(34, 11).Value = 4179.5  # K34 was 4519.4287
$ws.Cells.Item(34, 12).Value = 6553.5  # L34 was 6469.778
$ws.Cells.Item(34, 13).Value = -3977.5  # M34 was -4317.4287
$ws.Cells.Item(34, 14).Value = -6957.5  # N34 was -6873.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1004.9091  # H107 was 885.75
$ws.Cells.Item(107, 9).Value = 1037  # I107 was 937
$ws.Cells.Item(107, 10).Value = 992.875  # J107 was 868.6667
$ws.Cells.Item(107, 11).Value = 1037  # K107 was 937
$ws.Cells.Item(107, 12).Value = 992.875  # L107 was 868.6667
$ws.Cells.Item(107, 13).Value = 883  # M107 was 983
$ws.Cells.Item(107, 14).Value = -4832.875  # N107 was -4708.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 928.8889  # H34 was 1012.1429
$ws.Cells.Item(34, 9).Value = 239.8  # I34 was 224.75
$ws.Cells.Item(34, 10).Value = 1040.0322  # J34 was 1143.375
$ws.Cells.Item(34, 11).Value = 719.4000000000001  # K34 was 674.25
$ws.Cells.Item(34, 12).Value = 3120.0966  # L34 was 3430.125
$ws.Cells.Item(34, 13).Value = -635.4000000000001  # M34 was -590.25
$ws.Cells.Item(34, 14).Value = -3288.0966  # N34 was -3598.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 2738.087  # H39 was 2710.4707
$ws.Cells.Item(39, 9).Value = 700  # I39 was 700
$ws.Cells.Item(39, 10).Value = 2932.1904  # J39 was 2978.5334
$ws.Cells.Item(39, 11).Value = 2100  # K39 was 2100
$ws.Cells.Item(39, 12).Value = 8796.5712  # L39 was 8935.600199999999
$ws.Cells.Item(39, 13).Value = -1806  # M39 was -1806
$ws.Cells.Item(39, 14).Value = -9384.5712  # N39 was -9523.600199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 3212.375  # H55 was 3199.8823
$ws.Cells.Item(55, 9).Value = 0  # I55 was 0
$ws.Cells.Item(55, 10).Value = 3212.375  # J55 was 3199.8823
$ws.Cells.Item(55, 11).Value = 0  # K55 was 0
$ws.Cells.Item(55, 12).Value = 9637.125  # L55 was 9599.6469
$ws.Cells.Item(55, 14).Value = -9991.125  # N55 was -9953.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118, 8).Value = 2165.1667  # H118 was 2282.4546
$ws.Cells.Item(118, 9).Value = 1150  # I118 was 1326.3334
$ws.Cells.Item(118, 10).Value = 2503.5557  # J118 was 2641
$ws.Cells.Item(118, 11).Value = 3450  # K118 was 3979.0002
$ws.Cells.Item(118, 12).Value = 7510.6671  # L118 was 7923
$ws.Cells.Item(118, 13).Value = -2207  # M118 was -2736.0002
$ws.Cells.Item(118, 14).Value = -9996.667099999999  # N118 was -10409

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1908.9697  # H132 was 1837.5
$ws.Cells.Item(132, 9).Value = 1428.909  # I132 was 1428.909
$ws.Cells.Item(132, 10).Value = 2149  # J132 was 2003.963
$ws.Cells.Item(132, 11).Value = 12860.181  # K132 was 12860.181
$ws.Cells.Item(132, 12).Value = 19341  # L132 was 18035.667
$ws.Cells.Item(132, 13).Value = -10330.181  # M132 was -10330.181
$ws.Cells.Item(132, 14).Value = -24401  # N132 was -23095.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 4760.8696  # H136 was 4220.154
$ws.Cells.Item(136, 9).Value = 1336.25  # I136 was 1370
$ws.Cells.Item(136, 10).Value = 6587.3335  # J136 was 4843.625
$ws.Cells.Item(136, 11).Value = 4008.75  # K136 was 4110
$ws.Cells.Item(136, 12).Value = 19762.0005  # L136 was 14530.875
$ws.Cells.Item(136, 13).Value = 1091.25  # M136 was 990
$ws.Cells.Item(136, 14).Value = -29962.0005  # N136 was -24730.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 24783.408  # H123 was 33303.535
$ws.Cells.Item(123, 9).Value = 0  # I123 was 0
$ws.Cells.Item(123, 10).Value = 24783.408  # J123 was 33303.535
$ws.Cells.Item(123, 11).Value = 0  # K123 was 0
$ws.Cells.Item(123, 12).Value = 24783.408  # L123 was 33303.535
$ws.Cells.Item(123, 14).Value = -29683.408  # N123 was -38203.535

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 6930.6665  # H100 was 6304.615
$ws.Cells.Item(100, 9).Value = 10493.333  # I100 was 9280
$ws.Cells.Item(100, 10).Value = 4555.5557  # J100 was 2833.3333
$ws.Cells.Item(100, 11).Value = 10493.333  # K100 was 9280
$ws.Cells.Item(100, 12).Value = 4555.5557  # L100 was 2833.3333
$ws.Cells.Item(100, 13).Value = -9952.333000000001  # M100 was -8739
$ws.Cells.Item(100, 14).Value = -5637.5557  # N100 was -3915.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(105, 8).Value = 40000  # H105 was 0
$ws.Cells.Item(105, 9).Value = 0  # I105 was 0
$ws.Cells.Item(105, 10).Value = 40000  # J105 was 0
$ws.Cells.Item(105, 11).Value = 0  # K105 was 0
$ws.Cells.Item(105, 12).Value = 40000  # L105 was 0
$ws.Cells.Item(105, 14).Value = -46988  # N105 was None

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 948.8333  # H96 was 897.6667
$ws.Cells.Item(96, 9).Value = 846.5  # I96 was 846.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 290.06668  # H107 was 197.5
$ws.Cells.Item(107, 9).Value = 179.25  # I107 was 197.14285
$ws.Cells.Item(107, 10).Value = 733.3333  # J107 was 200
$ws.Cells.Item(107, 11).Value = 537.75  # K107 was 591.4285500000001
$ws.Cells.Item(107, 12).Value = 2199.9999  # L107 was 600
$ws.Cells.Item(107, 13).Value = 1382.25  # M107 was 1328.57145
$ws.Cells.Item(107, 14).Value = -6039.9999  # N107 was -4440

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 65429  # H124 was 40343.2
$ws.Cells.Item(124, 9).Value = 0  # I124 was 0
$ws.Cells.Item(124, 10).Value = 65429  # J124 was 40343.2
$ws.Cells.Item(124, 11).Value = 0  # K124 was 0
$ws.Cells.Item(124, 12).Value = 65429  # L124 was 40343.2
$ws.Cells.Item(124, 14).Value = -75249  # N124 was -50163.2
